$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark from its current location (it sits
#    between "el" and " esta presente para arreglarlo." in GRUPO10).
#    It will be re-created later at the very end of the document.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Delete the "*Influencers:" paragraph and the paragraph right
#    after it ("AriGameplays, Ahrileth, Ariana Grande.") completely.
# ------------------------------------------------------------------
$pInfluencers = $null
$pNames = $null
foreach ($p in $d.Paragraphs) {
    if ($pInfluencers -ne $null -and $pNames -eq $null) {
        $pNames = $p
    }
    if ($p.Range.Text -like "*Influencers:*") {
        $pInfluencers = $p
    }
}

if ($pInfluencers -ne $null -and $pNames -ne $null) {
    $delRange = $d.Range($pInfluencers.Range.Start, $pNames.Range.End)
    $delRange.Delete()
}

# ------------------------------------------------------------------
# 3) Re-insert the "_GoBack" bookmark, collapsed, at the start of the
#    last (now empty, trailing) paragraph of the document.
# ------------------------------------------------------------------
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $pLast.Range
$r.Collapse(1)
$r.Bookmarks.Add("_GoBack")
